$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column N: "Is Show Corporate" header (sharedStrings gets a 14th entry,
# used range grows from A1:M1 to A1:N1, row spans updates to 1:14).
$ws.Range("N1").Value = "Is Show Corporate"

# Give the new column a sensible custom width (matches the author's
# post-entry "best fit" auto-sizing of column N).
$ws.Columns.Item(14).ColumnWidth = 15.5

# Author's cursor ended up on N3 after adding the header.
$ws.Range("N3").Select() | Out-Null
